$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.091.09"
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("D3").Value = "2.625.35"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'517.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "'141.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("D9").Value = "'6.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("D11").Value = "'0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "3.095.90"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "58.115.08"
$ws.Range("E14").Value = "  -4.07%  "
$ws.Range("D15").Value = "'20.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").Value = "'0.0000135"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "2.632.11"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'4.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'334.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.91%  "
$ws.Range("D20").Value = "'10.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'64.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'0.422"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'7.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.26%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0786"
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("D29").Value = "'6.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D31").Value = "'152.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "'18.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "'0.903"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.20%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.99%  "
$ws.Range("D37").Value = "'36.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "'0.849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("E39").Value = "  -5.73%  "
$ws.Range("D40").Value = "'3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").Value = "'0.0967"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "'268.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.48%  "
$ws.Range("D45").Value = "'19.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("D47").Value = "'0.0532"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "2.027.36"
$ws.Range("E48").Value = "  -5.63%  "
$ws.Range("D49").Value = "'4.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("D50").Value = "'0.0226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "'18.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.83%  "
